# Generate Report for Handback
# Updates the localization-status workbook to reflect a completed handback:
# - Status text changes from "Ready for handoff" to "Handed back: in sync with en-US"
# - "Latest Target File" / "Latest Handback File" / "Latest Handback DateTime" columns
#   are populated for the zh-cn and de-de handback rows, with a hyperlink added on the
#   new "Latest Target File" cell.
# - A few columns are widened to fit the new content.

$wb = $excel.ActiveWorkbook

$baseUrl = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/4a37f7a417e17d15a91c77482cf334dba1349a51/e2e/"
$statusText = "Handed back: in sync with en-US"

# ---------------------------------------------------------------------------
# Overview sheet: status text + wider status columns (E, F)
# ---------------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = $statusText
$wsOverview.Range("F2").Value = $statusText
$wsOverview.Range("E3").Value = $statusText
$wsOverview.Range("F3").Value = $statusText
$wsOverview.Columns.Item(5).ColumnWidth = 29.14
$wsOverview.Columns.Item(6).ColumnWidth = 29.14

# ---------------------------------------------------------------------------
# Per-language detail sheets (zh-cn, de-de)
# ---------------------------------------------------------------------------
$languages = @(
    @{ Sheet = "zh-cn"; Xlf = "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.zh-cn.xlf"; HandbackDate = "2016-08-18 18:38:42" },
    @{ Sheet = "de-de"; Xlf = "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.de-de.xlf"; HandbackDate = "2016-08-18 18:38:49" }
)

foreach ($lang in $languages) {
    $ws = $wb.Worksheets.Item($lang.Sheet)

    # Status column (C) for both data rows
    $ws.Range("C2").Value = $statusText
    $ws.Range("C3").Value = $statusText

    # Widen Status (C) and Latest Handback File (J) columns
    $ws.Columns.Item(3).ColumnWidth = 29.14
    $ws.Columns.Item(10).ColumnWidth = 39.14

    # Latest Target File (I) + Latest Handback File (J) + Latest Handback DateTime (K)
    $ws.Range("I2").Value = "a.md"
    $ws.Range("J2").Value = $lang.Xlf
    $ws.Range("K2").Value = $lang.HandbackDate

    $ws.Range("I3").Value = "a.md"
    $ws.Range("J3").Value = $lang.Xlf
    $ws.Range("K3").Value = $lang.HandbackDate

    # Rebuild hyperlinks in row order (A2, I2, A3, I3) so relationship ids line up
    # the way the handback report generator emits them.
    $ws.Hyperlinks.Delete()
    $ws.Hyperlinks.Add($ws.Range("A2"), $baseUrl + "a.md", "", "", "a.md")
    $ws.Range("A2").Style = "HyperLink"

    $ws.Hyperlinks.Add($ws.Range("I2"), $baseUrl + "a.md", "", "", "a.md")
    $ws.Range("I2").Style = "HyperLink"

    $ws.Hyperlinks.Add($ws.Range("A3"), $baseUrl + "b.md", "", "", "b.md")
    $ws.Range("A3").Style = "HyperLink"

    $ws.Hyperlinks.Add($ws.Range("I3"), $baseUrl + "a.md", "", "", "a.md")
    $ws.Range("I3").Style = "HyperLink"
}
